$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.521.57'
$ws.Range("E2").Value = '  +5.31%  '
$ws.Range("D3").Value = '2.255.79'
$ws.Range("E3").Value = '  +4.58%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '232.31'
$ws.Range("E5").Value = '  +1.95%  '
$ws.Range("E6").Value = '  +2.22%  '
$ws.Range("D7").Value = '64.25'
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.411'
$ws.Range("E9").Value = '  +4.13%  '
$ws.Range("D10").Value = '59.45'
$ws.Range("E10").Value = '  +2.70%  '
$ws.Range("D11").Value = '0.0902'
$ws.Range("E11").Value = '  +5.61%  '
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").Value = '2.586.47'
$ws.Range("E13").Value = '  +4.41%  '
$ws.Range("D14").Value = '16.37'
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("E16").Value = '  +3.12%  '
$ws.Range("E17").Value = '  +3.26%  '
$ws.Range("D18").Value = '2.261.20'
$ws.Range("E18").Value = '  +4.03%  '
$ws.Range("D19").Value = '41.358.22'
$ws.Range("E19").Value = '  +5.01%  '
$ws.Range("D20").Value = '73.91'
$ws.Range("E20").Value = '  +2.95%  '
$ws.Range("D21").Value = '0.0₃0918'
$ws.Range("E21").Value = '  +8.06%  '
$ws.Range("E22").Value = '  +1.99%  '
$ws.Range("D23").Value = '251.48'
$ws.Range("E23").Value = '  +9.43%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  +5.62%  '
$ws.Range("D26").Value = '2.39'
$ws.Range("E26").Value = '  +1.81%  '
$ws.Range("D27").Value = '9.86'
$ws.Range("E27").Value = '  +3.73%  '
$ws.Range("D28").Value = '173.61'
$ws.Range("E28").Value = '  +0.82%  '
$ws.Range("D29").Value = '0.146'
$ws.Range("E29").Value = '  +2.52%  '
$ws.Range("D30").Value = '20.50'
$ws.Range("E30").Value = '  +3.32%  '
$ws.Range("D31").Value = '1.47'
$ws.Range("E31").Value = '  +3.93%  '
$ws.Range("E32").Value = '  +9.24%  '
$ws.Range("E33").Value = '  +2.68%  '
$ws.Range("E34").Value = '  +3.88%  '
$ws.Range("D35").Value = '5.03'
$ws.Range("E35").Value = '  +6.35%  '
$ws.Range("D36").Value = '0.0638'
$ws.Range("E36").Value = '  +3.59%  '
$ws.Range("D37").Value = '7.09'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").Value = '3.88'
$ws.Range("E38").Value = '  +9.19%  '
$ws.Range("E39").Value = '  +1.79%  '
$ws.Range("D40").Value = '0.000258'
$ws.Range("E40").Value = '  +63.23%  '
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("E42").Value = '  +4.22%  '
$ws.Range("D43").Value = '4.87'
$ws.Range("E43").Value = '  +12.22%  '
$ws.Range("D44").Value = '8.89'
$ws.Range("E44").Value = '  +14.22%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '18.00'
$ws.Range("E45").Value = '  +2.47%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '103.04'
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("E47").Value = '  +4.49%  '
$ws.Range("D48").Value = '1.513.28'
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("D49").Value = '0.0946'
$ws.Range("E49").Value = '  +1.83%  '
$ws.Range("E50").Value = '  +2.65%  '
$ws.Range("E51").Value = '  -0.85%  '
